$d = $word.ActiveDocument
$wdReplaceAll = 2
$wdFindWrap = 1

$d.Content.Find.Execute("19+77=96", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "7+6=13", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("81-67=14", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "9+86=95", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("2+2=4", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "49+42=91", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("27+32=59", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "40+46=86", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("44+2=46", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "58-0=58", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("89-79=10", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "29+41=70", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("62+36=98", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "54+39=93", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("69-17=52", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "91-79=12", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("65-65=0", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "53-48=5", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("96-89=7", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "49+20=69", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("98-90=8", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "98-67=31", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("1+87=88", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "9+15=24", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("74-57=17", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "64+19=83", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("51-8=43", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "30+61=91", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("80-59=21", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "28+44=72", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("82-45=37", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "78-52=26", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("48+16=64", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "79-58=21", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("37+22=59", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "91-57=34", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("32+44=76", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "57-25=32", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("27+38=65", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "49-40=9", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("24+25=49", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "62-27=35", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("98-78=20", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "39+3=42", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("56-5=51", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "2+28=30", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("26+51=77", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "41+44=85", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("48+20=68", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "84-7=77", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("86-7=79", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "35+1=36", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("91-67=24", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "66+1=67", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("44+54=98", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "47-22=25", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("64-27=37", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "22+9=31", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("46+18=64", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "59-23=36", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("51+5=56", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "71+22=93", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("40-26=14", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "40+16=56", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("26+50=76", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "3+43=46", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("0+60=60", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "26+37=63", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("44-28=16", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "66+6=72", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("55-27=28", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "27+37=64", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("39-19=20", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "51+35=86", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("76+14=90", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "67-12=55", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("78-53=25", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "9+54=63", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("56-17=39", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "33-17=16", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("66-56=10", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "70-26=44", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("62-18=44", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "1+49=50", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("31+17=48", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "22+26=48", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("52-6=46", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "7+8=15", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("79-29=50", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "71-11=60", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("75+10=85", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "39+48=87", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("68-37=31", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "99-54=45", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("42+13=55", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "15+63=78", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("62-23=39", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "82+17=99", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("40-24=16", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "3+55=58", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("17+17=34", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "55-6=49", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("5+21=26", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "20+33=53", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("54+0=54", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "38+20=58", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("22-4=18", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "89-39=50", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("91-7=84", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "38-13=25", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("26-3=23", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "59+39=98", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("61-6=55", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "55+29=84", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("28+22=50", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "55+44=99", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("68-62=6", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "9+2=11", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("54-2=52", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "19+60=79", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("56+36=92", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "3+16=19", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("82-36=46", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "10+1=11", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("35-3=32", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "74-0=74", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("13+33=46", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "22+55=77", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("83-13=70", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "9+27=36", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("62+26=88", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "87-22=65", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("37+51=88", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "86-11=75", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("99-63=36", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "70-60=10", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("26-8=18", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "76-12=64", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("42+44=86", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "0+79=79", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("97-36=61", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "60-13=47", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("46+48=94", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "13+67=80", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("96-93=3", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "91-21=70", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("50+7=57", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "64-6=58", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("10+76=86", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "99-1=98", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("36-29=7", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "72-34=38", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("78-30=48", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "40-33=7", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("24+61=85", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "65-28=37", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("20+49=69", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "37+44=81", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("16+24=40", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "40-1=39", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("65+27=92", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "15+58=73", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("1+4=5", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "8+19=27", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("24+59=83", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "19+60=79", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("49+43=92", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "51+26=77", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("94-35=59", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "24+35=59", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("0+96=96", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "75-49=26", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("63+0=63", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "87-74=13", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("35-19=16", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "16+1=17", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("83-74=9", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "29+38=67", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("84-17=67", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "32+10=42", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("10+67=77", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "82-44=38", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("81-47=34", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "14+67=81", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("63+22=85", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "22+64=86", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("42+53=95", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "57-54=3", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("17+65=82", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "59-2=57", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("4+60=64", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "99-50=49", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("39-11=28", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "10+6=16", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("88-12=76", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "69-47=22", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("82-0=82", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "8+13=21", $wdReplaceAll) | Out-Null
$d.Content.Find.Execute("28+7=35", $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, "58+14=72", $wdReplaceAll) | Out-Null
